# Insert a new weekly price record for "Comercializadora del Agro de
# Limarí - Poroto verde" right above the existing row 67, pushing the
# rest of the table (old rows 67-102) down by one row (new rows 68-103).
#
# The new record re-uses the same Mercado/Región/Producto/Unidad/Origen
# boilerplate as its neighbours and carries the latest weekly price
# observation (fecha 44489 = 2021-10-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 67..102 down to 68..103, leaving a blank row 67 to fill in.
$ws.Rows.Item(67).Insert()

$ws.Cells.Item(67, 1).Value = 2
$ws.Cells.Item(67, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(67, 3).Value = "Coquimbo"
$ws.Cells.Item(67, 4).Value = 44489
$ws.Cells.Item(67, 5).Value = 4
$ws.Cells.Item(67, 6).Value = 100112031
$ws.Cells.Item(67, 7).Value = "Poroto verde"
$ws.Cells.Item(67, 8).Value = "Magnum"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 800
$ws.Cells.Item(67, 11).Value = 43000
$ws.Cells.Item(67, 12).Value = 45000
$ws.Cells.Item(67, 13).Value = 44000
$ws.Cells.Item(67, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(67, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(67, 16).Value = 1760
$ws.Cells.Item(67, 17).Value = 25
$ws.Cells.Item(67, 18).Value = "Hortaliza"
